$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: update the two "1-20%" wordings to "1-40%" (creates new shared strings 17 & 18) ---
$ws1.Range("D3").Value = "2022-2032 varying from 1-40% of the population"
$ws1.Range("D4").Value = "Only in years between campaigns, varying from 1-40% of the population"

# --- Sheet1: row heights ---
$ws1.Rows.Item(1).RowHeight = 60
$ws1.Rows.Item(2).RowHeight = 60
$ws1.Rows.Item(3).RowHeight = 60
$ws1.Rows.Item(4).RowHeight = 60
$ws1.Rows.Item(5).RowHeight = 60
$ws1.Rows.Item(6).RowHeight = 79

# --- Sheet1: column widths (COM ColumnWidth = XML width - 0.8333333333333334) ---
$ws1.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws1.Columns.Item(2).ColumnWidth = 25.666666666666668
$ws1.Columns.Item(3).ColumnWidth = 27.166666666666668
$ws1.Columns.Item(4).ColumnWidth = 29.330729166666668

# --- Add the new worksheet "pop times" right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "pop times"

# Copy Sheet1's table (values + styles) into the new sheet as the starting point
$ws1.Range("A1:D6").Copy($ws2.Range("A1"))

# --- Sheet2: row heights (same pattern as sheet1, plus row 7 for totals) ---
$ws2.Rows.Item(1).RowHeight = 60
$ws2.Rows.Item(2).RowHeight = 60
$ws2.Rows.Item(3).RowHeight = 60
$ws2.Rows.Item(4).RowHeight = 60
$ws2.Rows.Item(5).RowHeight = 60
$ws2.Rows.Item(6).RowHeight = 79
$ws2.Rows.Item(7).RowHeight = 17

# --- Sheet2: column widths, matching Sheet1 ---
$ws2.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 27.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 29.330729166666668

# --- Sheet2: content edits, in the exact order the new text first appears ---
$ws2.Range("C2").Value = "2020-2035, varying from population x 5% to population x 7%"
$ws2.Range("C3").Value = "2021-2035 using population x 6%"
$ws2.Range("C4").Value = "2020-2035 using population x 6%"
$ws2.Range("E1").Value = "Number of different runs per scenario"
$ws2.Range("D7").Value = "Total iterations"
$ws2.Range("D4").Value = "Only in years between campaigns, varying the CD quantifer from population x 0% to population x 40%"
$ws2.Range("D3").Value = "2022-2032 varying the CD quantifer from population x 0% to population x 50%"

$ws2.Range("C5").Value = "2020-2035 using population x 6%"
$ws2.Range("C6").Value = "2020-2035 using population x 6%"

# --- Sheet2: numeric "runs per scenario" column ---
$ws2.Range("E2").Value = 3
$ws2.Range("E3").Value = 51
$ws2.Range("E4").Value = 41
$ws2.Range("E5").Value = 11
$ws2.Range("E6").Value = 11
$ws2.Range("E7").Formula = "=SUM(E2:E6)"

# --- Match formatting (style s=2 body style) on the new D7 "Total iterations" cell ---
$ws1.Range("C2").Copy($ws2.Range("D7"))
$ws2.Range("D7").Value = "Total iterations"

# --- Selections: Sheet1 keeps C2 selected (not the active tab); Sheet2 is the active tab at F4 ---
$ws1.Range("C2").Select()
$ws2.Activate()
$ws2.Range("F4").Select()
